$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 3
$ws.Range("D3").Value = "-"
$ws.Range("F3").Value = "['MCT-2A-Programação de Computadores', -, -, -]"

# Row 4
$ws.Range("D4").Value = "-"
$ws.Range("F4").Value = "['MCT-2A-Programação de Computadores', -, -, -]"

# Row 6
$ws.Range("D6").Value = "-"
$ws.Range("F6").Value = "['MCT-2A-Programação de Computadores', -, -, -]"

# Row 7
$ws.Range("D7").Value = "-"
$ws.Range("F7").Value = "['MCT-2A-Programação de Computadores', -, -, -]"

# Row 18
$ws.Range("B18").Value = "['ELM-1NA-Acionamentos Elétricos', -, 'ELM-1NA-Lógica de Programação', -]"
$ws.Range("C18").Value = "-"
$ws.Range("F18").Value = "['MEC-1NB-Comandos Eletricos', -, -, -]"

# Row 19
$ws.Range("B19").Value = "['ELM-1NA-Acionamentos Elétricos', -, 'ELM-1NA-Lógica de Programação', -]"
$ws.Range("C19").Value = "-"
$ws.Range("F19").Value = "['MEC-1NB-Comandos Eletricos', -, -, -]"

# Row 20
$ws.Range("B20").Value = "['ELM-1NA-Acionamentos Elétricos', -, 'ELM-1NA-Lógica de Programação', -]"
$ws.Range("C20").Value = "-"
$ws.Range("F20").Value = "['MEC-1NB-Comandos Eletricos', -, -, -]"

# Row 21
$ws.Range("B21").Value = "['ELM-1NA-Acionamentos Elétricos', -, 'ELM-1NA-Lógica de Programação', -]"
$ws.Range("C21").Value = "-"
$ws.Range("D21").Value = "[-, -, 'MEC-1NB-Comandos Eletricos', -]"
$ws.Range("F21").Value = "-"
